# Estrategia.docx — add the new "Migración de la tabla Maestra" section
# (a Título-styled heading paragraph followed by an explanatory paragraph)
# right after the existing final paragraph ("... con las características
# preestablecidas."), before the section properties.

$d = $word.ActiveDocument

# Anchor an empty (collapsed) range at the very end of the document's
# main story so the new content is appended after the last paragraph
# without touching/replacing any existing text.
$endPos = $d.Content.End
$r = $d.Range($endPos, $endPos)

$newSectionXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo"/><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Migración de la tabla Maestra </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:tab/><w:t xml:space="preserve">Al migrar los datos de los clientes de la tabla </w:t></w:r><w:r><w:rPr><w:i/><w:lang w:val="es-ES"/></w:rPr><w:t>Maestra</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> hacia la tabla </w:t></w:r><w:r><w:rPr><w:i/><w:lang w:val="es-ES"/></w:rPr><w:t>Clientes</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">, existe la nulidad en el campo </w:t></w:r><w:r><w:rPr><w:i/><w:lang w:val="es-ES"/></w:rPr><w:t>Provincia</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>, por lo cual se carga por defecto con el atributo “</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="es-ES"/></w:rPr><w:t>Migrada</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>”. De esta forma, dejamos en claro, que este cliente fue introducido al sistema a través de la migración y no contaba con una provincia.</w:t></w:r></w:p>'

$r.InsertXML($newSectionXml)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
